$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:H2").ClearContents()
$ws.Range("A1").Value = "No"
$ws.Range("B1").Value = "Sub Budget Code"
$ws.Range("C1").Value = "Sub Budget Name"
$ws.Range("D1").Value = "Work Code"
$ws.Range("E1").Value = "Work Name"
$ws.Range("F1").Value = "Product Code"
$ws.Range("G1").Value = "Product Name"
$ws.Range("H1").Value = "Currency Code"
$ws.Range("I1").Value = "Currency Name"
